$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497000000001
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("M2").Value = 0.3883076666666667
$ws.Range("N2").Value = 1.164923
$ws.Range("O2").Value = 0.1188638477168776
$ws.Range("P2").Value = 0.1188638477168776
$ws.Range("Q2").Value = 1.194110404636778
$ws.Range("R2").Value = 10.746993641731
$ws.Range("S2").Value = 0.003139519499501544
$ws.Range("T2").Value = 0.003139519499501544
$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497000000001
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.6829215134520935
$ws.Range("P3").Value = 0.6829215134520935
$ws.Range("Q3").Value = 6.860653600124444
$ws.Range("R3").Value = 61.74588240112
$ws.Range("S3").Value = 0.018037826044626
$ws.Range("T3").Value = 0.018037826044626
$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497000000001
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("O4").Value = 0.1982146388310289
$ws.Range("P4").Value = 0.1982146388310289
$ws.Range("Q4").Value = 1.991271249633667
$ws.Range("R4").Value = 17.921441246703
$ws.Range("S4").Value = 0.005235391043195301
$ws.Range("T4").Value = 0.005235391043195301
$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("M5").Value = 0.3883076666666667
$ws.Range("N5").Value = 1.164923
$ws.Range("O5").Value = 0.1188638477168776
$ws.Range("P5").Value = 0.1188638477168776
$ws.Range("Q5").Value = 24.839190654569
$ws.Range("R5").Value = 223.552715891121
$ws.Range("S5").Value = 0.065306460030031
$ws.Range("T5").Value = 0.06530646003003102
$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.6829215134520935
$ws.Range("P6").Value = 0.6829215134520935
$ws.Range("S6").Value = 0.3752123743136639
$ws.Range("T6").Value = 0.375212374313664
$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("O7").Value = 0.1982146388310289
$ws.Range("P7").Value = 0.1982146388310289
$ws.Range("S7").Value = 0.1089035618215781
$ws.Range("T7").Value = 0.1089035618215781
$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("M8").Value = 0.3883076666666667
$ws.Range("N8").Value = 1.164923
$ws.Range("O8").Value = 0.1188638477168776
$ws.Range("P8").Value = 0.1188638477168776
$ws.Range("Q8").Value = 19.17634242809222
$ws.Range("R8").Value = 172.58708185283
$ws.Range("S8").Value = 0.05041786818734504
$ws.Range("T8").Value = 0.05041786818734505
$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.6829215134520935
$ws.Range("P9").Value = 0.6829215134520935
$ws.Range("S9").Value = 0.2896713130938035
$ws.Range("T9").Value = 0.2896713130938036
$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("O10").Value = 0.1982146388310289
$ws.Range("P10").Value = 0.1982146388310289
$ws.Range("S10").Value = 0.08407568596625553
$ws.Range("T10").Value = 0.08407568596625555